{"js": "// Replace the multiplication problems in the practice table with new values.\n// Each mapping is a unique old->new string, applied via find/replace so the\n// run formatting (font, size, etc.) stays intact.\nconst replacements = [\n  [\"985\u00d72=\", \"803\u00d73=\"],\n  [\"576\u00d75=\", \"129\u00d78=\"],\n  [\"651\u00d73=\", \"510\u00d75=\"],\n  [\"523\u00d77=\", \"738\u00d76=\"],\n  [\"119\u00d78=\", \"256\u00d73=\"],\n  [\"633\u00d75=\", \"605\u00d75=\"],\n  [\"397\u00d78=\", \"152\u00d73=\"],\n  [\"558\u00d75=\", \"390\u00d77=\"],\n  [\"541\u00d72=\", \"857\u00d76=\"],\n  [\"813\u00d73=\", \"483\u00d75=\"],\n  [\"586\u00d75=\", \"334\u00d76=\"],\n  [\"637\u00d78=\", \"285\u00d74=\"],\n  [\"408\u00d76=\", \"645\u00d77=\"],\n  [\"983\u00d75=\", \"103\u00d79=\"],\n  [\"149\u00d77=\", \"938\u00d78=\"],\n  [\"763\u00d76=\", \"254\u00d75=\"],\n  [\"148\u00d74=\", \"663\u00d75=\"],\n  [\"713\u00d73=\", \"551\u00d77=\"],\n  [\"872\u00d77=\", \"214\u00d75=\"],\n  [\"493\u00d79=\", \"804\u00d73=\"],\n  [\"674\u00d74=\", \"940\u00d79=\"],\n  [\"735\u00d78=\", \"273\u00d74=\"],\n  [\"759\u00d78=\", \"404\u00d79=\"],\n  [\"720\u00d73=\", \"413\u00d76=\"],\n  [\"246\u00d77=\", \"523\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the multiplication problems in the practice table with new values.\n# Uses Word's Find/Replace on the whole document content so run formatting\n# (font, size, etc.) is preserved -- only the visible text token changes.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"985\u00d72=\"; New = \"803\u00d73=\"},\n    @{Old = \"576\u00d75=\"; New = \"129\u00d78=\"},\n    @{Old = \"651\u00d73=\"; New = \"510\u00d75=\"},\n    @{Old = \"523\u00d77=\"; New = \"738\u00d76=\"},\n    @{Old = \"119\u00d78=\"; New = \"256\u00d73=\"},\n    @{Old = \"633\u00d75=\"; New = \"605\u00d75=\"},\n    @{Old = \"397\u00d78=\"; New = \"152\u00d73=\"},\n    @{Old = \"558\u00d75=\"; New = \"390\u00d77=\"},\n    @{Old = \"541\u00d72=\"; New = \"857\u00d76=\"},\n    @{Old = \"813\u00d73=\"; New = \"483\u00d75=\"},\n    @{Old = \"586\u00d75=\"; New = \"334\u00d76=\"},\n    @{Old = \"637\u00d78=\"; New = \"285\u00d74=\"},\n    @{Old = \"408\u00d76=\"; New = \"645\u00d77=\"},\n    @{Old = \"983\u00d75=\"; New = \"103\u00d79=\"},\n    @{Old = \"149\u00d77=\"; New = \"938\u00d78=\"},\n    @{Old = \"763\u00d76=\"; New = \"254\u00d75=\"},\n    @{Old = \"148\u00d74=\"; New = \"663\u00d75=\"},\n    @{Old = \"713\u00d73=\"; New = \"551\u00d77=\"},\n    @{Old = \"872\u00d77=\"; New = \"214\u00d75=\"},\n    @{Old = \"493\u00d79=\"; New = \"804\u00d73=\"},\n    @{Old = \"674\u00d74=\"; New = \"940\u00d79=\"},\n    @{Old = \"735\u00d78=\"; New = \"273\u00d74=\"},\n    @{Old = \"759\u00d78=\"; New = \"404\u00d79=\"},\n    @{Old = \"720\u00d73=\"; New = \"413\u00d76=\"},\n    @{Old = \"246\u00d77=\"; New = \"523\u00d79=\"}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $r.New, 2)\n}\n"}
